$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "nombre" column (E)
$ws.Range("E1").Value = "nombre"
$ws.Range("E5").Value = "Francisco"

# Update turno for the appointment in row 5
$ws.Range("A5").Value = 5493804401611
$ws.Range("B5").Value = 45819

# Border around the new header cell (left+right thin)
$ws.Range("E1").Borders(7).LineStyle = 1
$ws.Range("E1").Borders(10).LineStyle = 1

$ws.Range("E5").Select()
